$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the bordered/centered formatting used by column A data cells (e.g. row 15)
# onto the two brand-new rows (16 and 17) before writing their values.
$ws.Cells.Item(15,1).Copy($ws.Cells.Item(16,1))
$ws.Cells.Item(15,1).Copy($ws.Cells.Item(17,1))

# Final contents for every data row (2-17): A, B(name), C, D, E(in_service)
$rows = @(
    @{ Row=2;  A=0;  B="line1"; C=7;  D=9;  E=$true  },
    @{ Row=3;  A=1;  B="line2"; C=9;  D=8;  E=$true  },
    @{ Row=4;  A=2;  B="line3"; C=8;  D=10; E=$true  },
    @{ Row=5;  A=3;  B="line4"; C=8;  D=11; E=$true  },
    @{ Row=6;  A=4;  B="line5"; C=10; D=5;  E=$true  },
    @{ Row=7;  A=5;  B="line6"; C=12; D=8;  E=$true  },
    @{ Row=8;  A=6;  B="line7"; C=14; D=11; E=$true  },
    @{ Row=9;  A=7;  B="line8"; C=16; D=9;  E=$false },
    @{ Row=10; A=8;  B="extr1"; C=5;  D=12; E=$true  },
    @{ Row=11; A=9;  B="extr2"; C=5;  D=9;  E=$true  },
    @{ Row=12; A=10; B="extr3"; C=10; D=11; E=$false },
    @{ Row=13; A=11; B="extr4"; C=7;  D=8;  E=$false },
    @{ Row=14; A=12; B="extr5"; C=9;  D=11; E=$true  },
    @{ Row=15; A=13; B="extr6"; C=7;  D=11; E=$false },
    @{ Row=16; A=14; B="extr7"; C=5;  D=7;  E=$true  },
    @{ Row=17; A=15; B="extr8"; C=8;  D=5;  E=$false }
)

foreach ($r in $rows) {
    $ws.Cells.Item($r.Row, 1).Value = $r.A
    $ws.Cells.Item($r.Row, 2).Value = $r.B
    $ws.Cells.Item($r.Row, 3).Value = $r.C
    $ws.Cells.Item($r.Row, 4).Value = $r.D
    $ws.Cells.Item($r.Row, 5).Value = $r.E
}

"done"
